$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear extra cells in row 4 (I4, K4:R4) that should no longer exist
$ws.Range("I4").Value = ""
$ws.Range("K4:R4").Value = ""

# Row 5
$ws.Range("A5").Value = "6VA35665"
$ws.Range("B5").Value = "B.BELLIES CLIP DUNCAN BEAR 8,5CM"
$ws.Range("C5").Value = "Consumo"
$ws.Range("D5").Value = "Tiene PT"
$ws.Range("E5").Value = "Tiene ES"
$ws.Range("F5").Value = "Tiene IT"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "1"
$ws.Range("H5").Value = "UND"
$ws.Range("J5").Value = "Solo Revisión"

# Row 6
$ws.Range("A6").Value = "6VA27994"
$ws.Range("B6").Value = "INVISIBOBBLE SPRUNCHIE DUO BRITISH ROYAL LADIES"
$ws.Range("C6").Value = "Consumo"
$ws.Range("D6").Value = "Tiene PT"
$ws.Range("E6").Value = "Tiene ES"
$ws.Range("F6").Value = "Tiene IT"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "2"
$ws.Range("H6").Value = "UND"
$ws.Range("J6").Value = "Solo Revisión"

# Row 7
$ws.Range("A7").Value = "6VA32949"
$ws.Range("B7").Value = "INVISIBOBBLE SPRUNCHIE DUO ITS SWEATER TIME"
$ws.Range("C7").Value = "Consumo"
$ws.Range("D7").Value = "Tiene PT"
$ws.Range("E7").Value = "Tiene ES"
$ws.Range("F7").Value = "Tiene IT"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "2"
$ws.Range("H7").Value = "UND"
$ws.Range("J7").Value = "Solo Revisión"

# Row 8
$ws.Range("A8").Value = "6VA14310"
$ws.Range("B8").Value = "GUYLOND LIMA UÑAS GRANDES 4 DISEÑOS"
$ws.Range("C8").Value = "Consumo"
$ws.Range("D8").Value = "Tiene PT"
$ws.Range("E8").Value = "Tiene ES"
$ws.Range("F8").Value = "Tiene IT"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4"
$ws.Range("H8").Value = "UND"
$ws.Range("J8").Value = "Solo Revisión"

# Row 9
$ws.Range("A9").Value = "6VA14310"
$ws.Range("B9").Value = "GUYLOND LIMA UÑAS GRANDES 4 DISEÑOS"
$ws.Range("C9").Value = "Consumo"
$ws.Range("D9").Value = "Tiene PT"
$ws.Range("E9").Value = "Tiene ES"
$ws.Range("F9").Value = "Tiene IT"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "4"
$ws.Range("H9").Value = "UND"
$ws.Range("I9").NumberFormat = "@"
$ws.Range("J9").Value = "Solo Revisión"
$ws.Range("K9:R9").NumberFormat = "@"
